$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting rows 8..end down by one.
$ws.Rows("8:8").Insert(-4121)

# Copy the formatting from A7:F7 (same logical "boolean field" row) into A8:F8
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 8 with the "Force" field (matches rows 7/9 boolean layout)
$ws.Range("A8").Value = "Force"
$ws.Range("B8:F8").Value = $false
